$d = $word.ActiveDocument

# Replace all occurrences of the COD-TEMPERATURE value 0.214 -> 0.545
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$d.Content.Find.Execute("0.214", $false, $false, $false, $false, $false, $true, 1, $false, "0.545", 2)

# Replace all occurrences of the HAKE-TEMPERATURE value 0.643 -> 0.455
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$d.Content.Find.Execute("0.643", $false, $false, $false, $false, $false, $true, 1, $false, "0.455", 2)

$d.Save()
